$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 17: Inscritos 100 -> 101, Pagos 43 -> 44, Inscricoes homologadas 43 -> 44
$ws.Range("E17").Value = 101
$ws.Range("F17").Value = 44
$ws.Range("H17").Value = 44

# Row 19: Inscritos 48 -> 49
$ws.Range("E19").Value = 49

# Row 33: Inscritos 37 -> 38
$ws.Range("E33").Value = 38

# Row 38: Inscritos 65 -> 66
$ws.Range("E38").Value = 66

# Row 42: Inscritos 30 -> 31
$ws.Range("E42").Value = 31

# Row 49: Inscritos 58 -> 59
$ws.Range("E49").Value = 59
